# Applies updated market-board / profit figures to each class sheet,
# as refreshed by the scheduled pricing-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 7601.6665
$ws.Range("I43").Value = 10437.25
$ws.Range("K43").Value = 10437.25
$ws.Range("M43").Value = -10368.25
$ws.Range("H53").Value = 904.53845
$ws.Range("I53").Value = 96.42856999999999
$ws.Range("K53").Value = 96.42856999999999
$ws.Range("M53").Value = 540.57143
$ws.Range("H62").Value = 6751.36
$ws.Range("I62").Value = 6298.067
$ws.Range("K62").Value = 6298.067
$ws.Range("M62").Value = -5674.067
$ws.Range("H65").Value = 6751.36
$ws.Range("I65").Value = 6298.067
$ws.Range("K65").Value = 31490.335
$ws.Range("M65").Value = -28370.335
$ws.Range("H129").Value = 1094.0834
$ws.Range("I129").Value = 859.1111
$ws.Range("K129").Value = 2577.3333
$ws.Range("M129").Value = 2422.6667
$ws.Range("H137").Value = 8502.726000000001
$ws.Range("I137").Value = 14010.88
$ws.Range("J137").Value = 3206.423
$ws.Range("K137").Value = 42032.64
$ws.Range("L137").Value = 9619.269
$ws.Range("M137").Value = -39482.64
$ws.Range("N137").Value = -14719.269
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()
$ws.Range("H97").Value = 2737.5715
$ws.Range("I97").Value = 2350.4
$ws.Range("J97").Value = 3184.3076
$ws.Range("K97").Value = 2350.4
$ws.Range("L97").Value = 3184.3076
$ws.Range("M97").Value = -1854.4
$ws.Range("N97").Value = -4176.3076
$ws.Range("H101").Value = 53825.25
$ws.Range("J101").Value = 53825.25
$ws.Range("L101").Value = 53825.25
$ws.Range("N101").Value = -60315.25
$ws.Range("H109").Value = 106499.5
$ws.Range("J109").Value = 106499.5
$ws.Range("L109").Value = 106499.5
$ws.Range("N109").Value = -109273.5
$ws.Range("H110").Value = 1289.25
$ws.Range("I110").Value = 1298.25
$ws.Range("K110").Value = 1298.25
$ws.Range("M110").Value = 746.75
$ws.Range("H132").Value = 1020.6818
$ws.Range("I132").Value = 928.1395
$ws.Range("K132").Value = 2784.4185
$ws.Range("M132").Value = -254.4184999999998
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 25000
$ws.Range("J35").Value = 25000
$ws.Range("L35").Value = 25000
$ws.Range("N35").Value = -25620
$ws.Range("H75").Value = 26218.555
$ws.Range("J75").Value = 34995.5
$ws.Range("L75").Value = 34995.5
$ws.Range("N75").Value = -36867.5
$ws.Range("H78").Value = 26218.555
$ws.Range("J78").Value = 34995.5
$ws.Range("L78").Value = 104986.5
$ws.Range("N78").Value = -114346.5
$ws.Range("H86").Value = 1119.875
$ws.Range("I86").Value = 1116.1852
$ws.Range("K86").Value = 1116.1852
$ws.Range("M86").Value = 6.814800000000105
$ws.Range("H89").Value = 1119.875
$ws.Range("I89").Value = 1116.1852
$ws.Range("K89").Value = 5580.925999999999
$ws.Range("M89").Value = 35.07400000000052
$ws.Range("H107").Value = 22571
$ws.Range("I107").Value = 26850.25
$ws.Range("K107").Value = 26850.25
$ws.Range("M107").Value = -24930.25
$ws.Range("H134").Value = 1936.8276
$ws.Range("I134").Value = 1704.5283
$ws.Range("K134").Value = 5113.5849
$ws.Range("M134").Value = -2578.5849
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 305.36365
$ws.Range("I7").Value = 338.88235
$ws.Range("K7").Value = 338.88235
$ws.Range("M7").Value = -225.88235
$ws.Range("H41").Value = 13065.333
$ws.Range("I41").Value = 2133.3333
$ws.Range("K41").Value = 2133.3333
$ws.Range("M41").Value = -1705.3333
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H51").Value = 31858.1
$ws.Range("J51").Value = 31858.1
$ws.Range("L51").Value = 31858.1
$ws.Range("N51").Value = -33330.1
$ws.Range("H58").Value = 1716.4333
$ws.Range("I58").Value = 1603.2413
$ws.Range("J58").Value = 4999
$ws.Range("K58").Value = 1603.2413
$ws.Range("L58").Value = 4999
$ws.Range("M58").Value = -1400.2413
$ws.Range("N58").Value = -5405
$ws.Range("H59").Value = 47024.5
$ws.Range("J59").Value = 47024.5
$ws.Range("L59").Value = 47024.5
$ws.Range("N59").Value = -49314.5
$ws.Range("H61").Value = 31858.1
$ws.Range("J61").Value = 31858.1
$ws.Range("L61").Value = 31858.1
$ws.Range("N61").Value = -32554.1
$ws.Range("H98").Value = 62182.25
$ws.Range("J98").Value = 62182.25
$ws.Range("L98").Value = 62182.25
$ws.Range("N98").Value = -66674.25
$ws.Range("H136").Value = 1716.4333
$ws.Range("I136").Value = 1603.2413
$ws.Range("J136").Value = 4999
$ws.Range("K136").Value = 4809.7239
$ws.Range("L136").Value = 14997
$ws.Range("M136").Value = -2259.7239
$ws.Range("N136").Value = -20097
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 37414290
$ws.Range("I4").Value = 45855270
$ws.Range("K4").Value = 137565810
$ws.Range("M4").Value = -137565698
$ws.Range("H18").Value = 453.1
$ws.Range("I18").Value = 403.44446
$ws.Range("J18").Value = 900
$ws.Range("K18").Value = 1210.33338
$ws.Range("L18").Value = 2700
$ws.Range("M18").Value = -1041.33338
$ws.Range("N18").Value = -3038
$ws.Range("H117").Value = 2376.25
$ws.Range("J117").Value = 2460.75
$ws.Range("L117").Value = 7382.25
$ws.Range("N117").Value = -14266.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 52333
$ws.Range("J134").Value = 52333
$ws.Range("L134").Value = 156999
$ws.Range("N134").Value = -162069
$ws.Range("H136").Value = 7078
$ws.Range("J136").Value = 7078
$ws.Range("L136").Value = 21234
$ws.Range("N136").Value = -26334
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2674.111
$ws.Range("I16").Value = 2962.6667
$ws.Range("J16").Value = 2097
$ws.Range("K16").Value = 2962.6667
$ws.Range("L16").Value = 2097
$ws.Range("M16").Value = -2792.6667
$ws.Range("N16").Value = -2437
$ws.Range("H97").Value = 31672.5
$ws.Range("J97").Value = 31672.5
$ws.Range("L97").Value = 31672.5
$ws.Range("N97").Value = -33654.5
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H127").Value = 116249.5
$ws.Range("J127").Value = 116249.5
$ws.Range("L127").Value = 116249.5
$ws.Range("N127").Value = -126169.5
$ws.Range("H132").Value = 5289.7295
$ws.Range("I132").Value = 4529.871
$ws.Range("J132").Value = 9215.666999999999
$ws.Range("K132").Value = 13589.613
$ws.Range("L132").Value = 27647.001
$ws.Range("M132").Value = -11059.613
$ws.Range("N132").Value = -32707.001
$ws.Range("H136").Value = 5036.8
$ws.Range("I136").Value = 3242
$ws.Range("K136").Value = 9726
$ws.Range("M136").Value = -7176
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1491.4
$ws.Range("I2").Value = 666.6667
$ws.Range("J2").Value = 2728.5
$ws.Range("K2").Value = 666.6667
$ws.Range("L2").Value = 2728.5
$ws.Range("M2").Value = -554.6667
$ws.Range("N2").Value = -2952.5
$ws.Range("H41").Value = 19948
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 19948
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 19948
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -20728
$ws.Range("H107").Value = 1763.7
$ws.Range("I107").Value = 1763.7
$ws.Range("K107").Value = 5291.1
$ws.Range("M107").Value = -3371.1
$ws.Range("H113").Value = 1294.8462
$ws.Range("I113").Value = 860.8570999999999
$ws.Range("K113").Value = 2582.5713
$ws.Range("M113").Value = -412.5712999999996
$ws.Range("H136").Value = 29882.912
$ws.Range("I136").Value = 31234.863
$ws.Range("J136").Value = 140
$ws.Range("K136").Value = 93704.58900000001
$ws.Range("L136").Value = 420
$ws.Range("M136").Value = -91154.58900000001
$ws.Range("N136").Value = -5520
